{"js": "// Each table cell holds one two-digit-by-two-digit multiplication\n// prompt as a single text run, e.g. \"32\\u00d782=\". Replace each old\n// prompt with its new prompt, in document order, via search + Replace.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"32\u00d782=\", \"11\u00d790=\"],\n  [\"89\u00d782=\", \"67\u00d780=\"],\n  [\"24\u00d765=\", \"78\u00d785=\"],\n  [\"40\u00d777=\", \"51\u00d743=\"],\n  [\"87\u00d738=\", \"21\u00d758=\"],\n  [\"23\u00d755=\", \"93\u00d746=\"],\n  [\"40\u00d724=\", \"64\u00d757=\"],\n  [\"84\u00d736=\", \"92\u00d797=\"],\n  [\"49\u00d762=\", \"14\u00d715=\"],\n  [\"34\u00d761=\", \"64\u00d768=\"],\n  [\"61\u00d784=\", \"63\u00d752=\"],\n  [\"46\u00d743=\", \"33\u00d743=\"],\n  [\"18\u00d774=\", \"40\u00d725=\"],\n  [\"57\u00d742=\", \"85\u00d719=\"],\n  [\"39\u00d787=\", \"63\u00d711=\"],\n  [\"89\u00d788=\", \"30\u00d790=\"],\n  [\"72\u00d791=\", \"49\u00d733=\"],\n  [\"23\u00d771=\", \"34\u00d761=\"],\n  [\"33\u00d778=\", \"92\u00d765=\"],\n  [\"49\u00d716=\", \"78\u00d761=\"],\n  [\"67\u00d743=\", \"37\u00d741=\"],\n  [\"37\u00d757=\", \"56\u00d786=\"],\n  [\"60\u00d735=\", \"38\u00d713=\"],\n  [\"77\u00d724=\", \"45\u00d766=\"],\n  [\"28\u00d735=\", \"93\u00d723=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text: ${oldText}`);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Each table cell holds one two-digit-by-two-digit multiplication\n# prompt as a single text run, e.g. \"32\u00d782=\". Replace each old\n# prompt with its new prompt, in document order, via Find/Replace.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"32\u00d782=\", \"11\u00d790=\"),\n  @(\"89\u00d782=\", \"67\u00d780=\"),\n  @(\"24\u00d765=\", \"78\u00d785=\"),\n  @(\"40\u00d777=\", \"51\u00d743=\"),\n  @(\"87\u00d738=\", \"21\u00d758=\"),\n  @(\"23\u00d755=\", \"93\u00d746=\"),\n  @(\"40\u00d724=\", \"64\u00d757=\"),\n  @(\"84\u00d736=\", \"92\u00d797=\"),\n  @(\"49\u00d762=\", \"14\u00d715=\"),\n  @(\"34\u00d761=\", \"64\u00d768=\"),\n  @(\"61\u00d784=\", \"63\u00d752=\"),\n  @(\"46\u00d743=\", \"33\u00d743=\"),\n  @(\"18\u00d774=\", \"40\u00d725=\"),\n  @(\"57\u00d742=\", \"85\u00d719=\"),\n  @(\"39\u00d787=\", \"63\u00d711=\"),\n  @(\"89\u00d788=\", \"30\u00d790=\"),\n  @(\"72\u00d791=\", \"49\u00d733=\"),\n  @(\"23\u00d771=\", \"34\u00d761=\"),\n  @(\"33\u00d778=\", \"92\u00d765=\"),\n  @(\"49\u00d716=\", \"78\u00d761=\"),\n  @(\"67\u00d743=\", \"37\u00d741=\"),\n  @(\"37\u00d757=\", \"56\u00d786=\"),\n  @(\"60\u00d735=\", \"38\u00d713=\"),\n  @(\"77\u00d724=\", \"45\u00d766=\"),\n  @(\"28\u00d735=\", \"93\u00d723=\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Could not find text: $oldText\"\n  }\n}\n"}
